$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: B2/C2 = Status ("In Translation" -> "Ready for handoff"),
# D2 = Latest Handoff Date ("2016-27-19 20:27:35" -> "2016-28-19 20:28:21")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-28-19 20:28:21"

# zh-cn sheet: C2 = Status, E2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-19 20:28:19"

# de-de sheet: C2 = Status, E2 = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-19 20:28:21"
